{"js": "// Update in-text citation placeholders with the checked citation keys\n// (per commit: \"update: datasets | citation check | document browser\").\nconst body = context.document.body;\n\nconst replacements = [\n  { find: \"Ref-A1B2C3\", replace: \"Ref-s544112\" },\n  { find: \"Ref-J7X2B9\", replace: \"Ref-f465817\" },\n  { find: \"Ref-DJ49F2\", replace: \"Pearse et al., 2001\" },\n  { find: \"Ref-AB1CD2\", replace: \"Ref-f240694\" },\n  { find: \"Ref-DJ49KL\", replace: \"Ref-f214141\" }\n];\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update in-text citation placeholders with the checked citation keys\n# (per commit: \"update: datasets | citation check | document browser\").\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$pairs = @(\n    @{ Find = \"Ref-A1B2C3\"; Replace = \"Ref-s544112\" },\n    @{ Find = \"Ref-J7X2B9\"; Replace = \"Ref-f465817\" },\n    @{ Find = \"Ref-DJ49F2\"; Replace = \"Pearse et al., 2001\" },\n    @{ Find = \"Ref-AB1CD2\"; Replace = \"Ref-f240694\" },\n    @{ Find = \"Ref-DJ49KL\"; Replace = \"Ref-f214141\" }\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair.Find, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $pair.Replace, $wdReplaceAll)\n}\n"}
